$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 26.1813690087763
$ws.Range("G2").Value = 22.08919175257506
$ws.Range("H2").Value = 30.41838090909802
$ws.Range("I2").Value = 3.355703362776727
$ws.Range("J2").Value = 2.04995262557871
$ws.Range("K2").Value = 4.897767566540358
$ws.Range("L2").Value = 0.2452766279584661
$ws.Range("M2").Value = 0.1626529843809138
$ws.Range("N2").Value = 0.3486432408210731

# Row 3
$ws.Range("F3").Value = 0.9772969054597146
$ws.Range("G3").Value = 0.03730947449461339
$ws.Range("H3").Value = 1.999502663789804
$ws.Range("I3").Value = 0.8414889294528893
$ws.Range("J3").Value = 0.03453000811037173
$ws.Range("K3").Value = 1.690550276461255
$ws.Range("L3").Value = 0.9934997885860725
$ws.Range("M3").Value = 0.03765563306937082
$ws.Range("N3").Value = 2.030154005518863

# Row 4
$ws.Range("F4").Value = 27.15866591423602
$ws.Range("G4").Value = 22.12650122706968
$ws.Range("H4").Value = 32.41788357288782
$ws.Range("I4").Value = 4.197192292229616
$ws.Range("J4").Value = 2.084482633689082
$ws.Range("K4").Value = 6.588317843001612
$ws.Range("L4").Value = 1.238776416544539
$ws.Range("M4").Value = 0.2003086174502846
$ws.Range("N4").Value = 2.378797246339936
